$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.3064541150169984
$ws.Range("I2").Value = 0.3064541150169984
$ws.Range("L2").Value = 4.794854564667076
$ws.Range("M2").Value = "[-3.5183936022065057, 13.108102731540658]"
$ws.Range("N2").Value = 0.2514902790445142
$ws.Range("O2").Value = 0.2514902790445142
$ws.Range("P2").Value = -0.3522105877898465
$ws.Range("Q2").Value = "[-3.4780795544247356, 2.7736583788450426]"
$ws.Range("R2").Value = 0.8214972736118136
$ws.Range("S2").Value = 0.8214972736118136
$ws.Range("T2").Value = 11.66594461504863
$ws.Range("U2").Value = "[7.317204321260746, 16.014684908836507]"
$ws.Range("V2").Value = 0.000002378648665812833
$ws.Range("W2").Value = 0.000002378648665812833
$ws.Range("X2").Value = 1.456896896896932
$ws.Range("Y2").Value = -11.47306306306334
$ws.Range("Z2").Value = 14.3868568568572
$ws.Range("H3").Value = 0.2739869732706874
$ws.Range("I3").Value = 0.2739869732706874
$ws.Range("L3").Value = 4.561121824203502
$ws.Range("M3").Value = "[-3.4346357812810595, 12.556879429688063]"
$ws.Range("N3").Value = 0.256654998318695
$ws.Range("O3").Value = 0.256654998318695
$ws.Range("P3").Value = -0.8050527720910781
$ws.Range("Q3").Value = "[-3.94350068828989, 2.333395144107734]"
$ws.Range("R3").Value = 0.6079347573312086
$ws.Range("S3").Value = 0.6079347573312086
$ws.Range("T3").Value = 11.72891420521284
$ws.Range("U3").Value = "[7.6308027726797265, 15.827025637745958]"
$ws.Range("V3").Value = 0.0000006976820634552183
$ws.Range("W3").Value = 0.0000006976820634552183
$ws.Range("X3").Value = 3.330050050050129
$ws.Range("Y3").Value = -9.651941941942173
$ws.Range("Z3").Value = 16.31204204204243
$ws.Range("H4").Value = 0.02914682919290212
$ws.Range("I4").Value = 0.02914682919290212
$ws.Range("L4").Value = 9.985011143156427
$ws.Range("M4").Value = "[0.4318360758569497, 19.538186210455905]"
$ws.Range("N4").Value = 0.04089300938520535
$ws.Range("O4").Value = 0.04089300938520535
$ws.Range("P4").Value = -2.012631930227695
$ws.Range("Q4").Value = "[-3.3208426848756964, -0.704421175579693]"
$ws.Range("R4").Value = 0.003346328941514098
$ws.Range("S4").Value = 0.003346328941514098
$ws.Range("T4").Value = 14.90159061763793
$ws.Range("U4").Value = "[9.73175853562951, 20.07142269964636]"
$ws.Range("V4").Value = 0.0000006065859068371537
$ws.Range("W4").Value = 0.0000006065859068371537
$ws.Range("X4").Value = 8.325125125125322
$ws.Range("Y4").Value = 2.913793793793861
$ws.Range("Z4").Value = 13.73645645645678
$ws.Range("H5").Value = 0.4421755868177183
$ws.Range("I5").Value = 0.4421755868177183
$ws.Range("L5").Value = 5.026454883750497
$ws.Range("M5").Value = "[-4.482605682794055, 14.535515450295048]"
$ws.Range("N5").Value = 0.2927148762421075
$ws.Range("O5").Value = 0.2927148762421075
$ws.Range("P5").Value = -3.018947895341543
$ws.Range("Q5").Value = "[-6.151106336758394, 0.113210546075309]"
$ws.Range("R5").Value = 0.05849927387699227
$ws.Range("S5").Value = 0.05849927387699227
$ws.Range("T5").Value = 16.46864060885407
$ws.Range("U5").Value = "[11.132866431031271, 21.80441478667687]"
$ws.Range("V5").Value = 0.0000001488972116714393
$ws.Range("W5").Value = 0.0000001488972116714393
$ws.Range("X5").Value = 12.48768768768799
$ws.Range("Y5").Value = -0.4682882882883028
$ws.Range("Z5").Value = 25.44366366366427
$ws.Range("B6").Value = 1
$ws.Range("H6").Value = 0.01330535525417875
$ws.Range("I6").Value = 0.01330535525417875
$ws.Range("L6").Value = 11.21573166963026
$ws.Range("M6").Value = "[2.786553110680213, 19.6449102285803]"
$ws.Range("N6").Value = 0.01025121893898451
$ws.Range("O6").Value = 0.01025121893898451
$ws.Range("P6").Value = 3.113290017070967
$ws.Range("Q6").Value = "[2.1069740519571187, 4.119605982184814]"
$ws.Range("R6").Value = 0.0000001415944284044457
$ws.Range("S6").Value = 0.0000001415944284044457
$ws.Range("T6").Value = 14.58529900922494
$ws.Range("U6").Value = "[9.572299107475468, 19.598298910974407]"
$ws.Range("V6").Value = 0.000000503613034608108
$ws.Range("W6").Value = 0.000000503613034608108
$ws.Range("X6").Value = 13.11207207207238
$ws.Range("Y6").Value = 8.949509509509721
$ws.Range("Z6").Value = 17.27463463463505
$ws.Range("H7").Value = 0.5477543247287908
$ws.Range("I7").Value = 0.5477543247287908
$ws.Range("L7").Value = 4.329255128358434
$ws.Range("M7").Value = "[-6.296864444416708, 14.955374701133575]"
$ws.Range("N7").Value = 0.4162120663314171
$ws.Range("O7").Value = 0.4162120663314171
$ws.Range("P7").Value = 2.647868883205812
$ws.Range("Q7").Value = "[-0.4842895582110387, 5.7800273246226626]"
$ws.Range("R7").Value = 0.09552586427299792
$ws.Range("S7").Value = 0.09552586427299792
$ws.Range("T7").Value = 16.46029286492948
$ws.Range("U7").Value = "[10.858549023020998, 22.062036706837965]"
$ws.Range("V7").Value = 0.0000004127864443592699
$ws.Range("W7").Value = 0.0000004127864443592699
$ws.Range("X7").Value = 15.03725725725762
$ws.Range("Y7").Value = 2.081281281281329
$ws.Range("Z7").Value = 27.9932332332339
$ws.Range("H8").Value = 0.4892309729125598
$ws.Range("I8").Value = 0.4892309729125598
$ws.Range("L8").Value = 4.145831155309981
$ws.Range("M8").Value = "[-5.478171930864882, 13.769834241484844]"
$ws.Range("N8").Value = 0.3901969981087108
$ws.Range("O8").Value = 0.3901969981087108
$ws.Range("P8").Value = 1.792500312859041
$ws.Range("Q8").Value = "[-1.3333686537758478, 4.91836927949393]"
$ws.Range("R8").Value = 0.2541995405117288
$ws.Range("S8").Value = 0.2541995405117288
$ws.Range("T8").Value = 13.17474453329475
$ws.Range("U8").Value = "[8.050588019023696, 18.29890104756581]"
$ws.Range("V8").Value = 0.000005066353445881688
$ws.Range("W8").Value = 0.000005066353445881688
$ws.Range("X8").Value = 18.57543543543588
$ws.Range("Y8").Value = 5.645475475475612
$ws.Range("Z8").Value = 31.50539539539615
$ws.Range("H9").Value = 0.7391345481616722
$ws.Range("I9").Value = 0.7391345481616722
$ws.Range("L9").Value = 2.877431347799724
$ws.Range("M9").Value = "[-7.583548802707168, 13.338411498306616]"
$ws.Range("N9").Value = 0.5823175798635276
$ws.Range("O9").Value = 0.5823175798635276
$ws.Range("P9").Value = 2.270500396288119
$ws.Range("Q9").Value = "[-0.8616580451287321, 5.40265883770497]"
$ws.Range("R9").Value = 0.1512304904986119
$ws.Range("S9").Value = 0.1512304904986119
$ws.Range("T9").Value = 15.30349701074203
$ws.Range("U9").Value = "[9.876994647525525, 20.72999937395853]"
$ws.Range("V9").Value = 0.0000009298065724028248
$ws.Range("W9").Value = 0.0000009298065724028248
$ws.Range("X9").Value = 16.59821821821862
$ws.Range("Y9").Value = 3.642242242242327
$ws.Range("Z9").Value = 29.55419419419491
$ws.Range("F10").Value = 24.07000000000032
$ws.Range("H10").Value = 0.1328851661852331
$ws.Range("I10").Value = 0.1328851661852331
$ws.Range("L10").Value = 6.13647690040034
$ws.Range("M10").Value = "[-1.7224861380455465, 13.995439938846227]"
$ws.Range("N10").Value = 0.1228023636215121
$ws.Range("O10").Value = 0.1228023636215121
$ws.Range("P10").Value = 1.767342413731195
$ws.Range("Q10").Value = "[0.006289474781961069, 3.5283953526804295]"
$ws.Range("R10").Value = 0.04922134428305269
$ws.Range("S10").Value = 0.04922134428305269
$ws.Range("T10").Value = 14.41440503242287
$ws.Range("U10").Value = "[10.159426159439386, 18.669383905406363]"
$ws.Range("V10").Value = 0.00000001865228882458325
$ws.Range("W10").Value = 0.00000001865228882458325
$ws.Range("X10").Value = 17.29955955955979
$ws.Range("Y10").Value = 10.55321321321335
$ws.Range("Z10").Value = 24.04590590590623
$ws.Range("F11").Value = 24.07000000000032
$ws.Range("H11").Value = 0.1986529681953104
$ws.Range("I11").Value = 0.1986529681953104
$ws.Range("L11").Value = 7.429759131598314
$ws.Range("M11").Value = "[-3.4523063811305468, 18.311824644327174]"
$ws.Range("N11").Value = 0.1758977890582818
$ws.Range("O11").Value = 0.1758977890582818
$ws.Range("P11").Value = 2.471763589310888
$ws.Range("Q11").Value = "[-0.6226580034141929, 5.566185182035969]"
$ws.Range("R11").Value = 0.1146484650294295
$ws.Range("S11").Value = 0.1146484650294295
$ws.Range("T11").Value = 13.10903740708836
$ws.Range("U11").Value = "[7.3304195646374595, 18.887655249539264]"
$ws.Range("V11").Value = 0.00003803715874384217
$ws.Range("W11").Value = 0.00003803715874384217
$ws.Range("X11").Value = 14.60102102102122
$ws.Range("Y11").Value = 2.746726726726765
$ws.Range("Z11").Value = 26.45531531531567
$ws.Range("F12").Value = 24.07000000000032
$ws.Range("H12").Value = 0.07191338858868879
$ws.Range("I12").Value = 0.07191338858868879
$ws.Range("L12").Value = 7.361025883856449
$ws.Range("M12").Value = "[-0.7785740985223732, 15.500625866235271]"
$ws.Range("N12").Value = 0.07519018355829599
$ws.Range("O12").Value = 0.07519018355829599
$ws.Range("P12").Value = 2.018921405009657
$ws.Range("Q12").Value = "[0.1823947686768843, 3.855448041342429]"
$ws.Range("R12").Value = 0.03192870084651522
$ws.Range("S12").Value = 0.03192870084651522
$ws.Range("T12").Value = 12.90873503223735
$ws.Range("U12").Value = "[8.468901406543612, 17.348568657931093]"
$ws.Range("V12").Value = 0.0000005106297109591651
$ws.Range("W12").Value = 0.0000005106297109591651
$ws.Range("X12").Value = 16.33579579579602
$ws.Range("Y12").Value = 9.300320320320445
$ws.Range("Z12").Value = 23.37127127127159
$ws.Range("F13").Value = 24.07000000000032
$ws.Range("H13").Value = 0.3962642393529212
$ws.Range("I13").Value = 0.3962642393529212
$ws.Range("L13").Value = 4.566007804513287
$ws.Range("M13").Value = "[-3.977602875387859, 13.109618484414433]"
$ws.Range("N13").Value = 0.2874819231889256
$ws.Range("O13").Value = 0.2874819231889256
$ws.Range("P13").Value = 1.83023716155081
$ws.Range("Q13").Value = "[-1.2956318050840787, 4.9561061281856995]"
$ws.Range("R13").Value = 0.2444826577796433
$ws.Range("S13").Value = 0.2444826577796433
$ws.Range("T13").Value = 11.95590019596423
$ws.Range("U13").Value = "[7.219439568581137, 16.69236082334732]"
$ws.Range("V13").Value = 0.000006949857148219607
$ws.Range("W13").Value = 0.000006949857148219607
$ws.Range("X13").Value = 17.05861861861885
$ws.Range("Y13").Value = 5.083853853853922
$ws.Range("Z13").Value = 29.03338338338377
$ws.Range("F14").Value = 24.07000000000032
$ws.Range("H14").Value = 0.5411232121008607
$ws.Range("I14").Value = 0.5411232121008607
$ws.Range("L14").Value = 3.569991751107207
$ws.Range("M14").Value = "[-4.8763862504004605, 12.016369752614874]"
$ws.Range("N14").Value = 0.3991147881845394
$ws.Range("O14").Value = 0.3991147881845394
$ws.Range("P14").Value = 1.239026532046425
$ws.Range("Q14").Value = "[-1.8868424345884645, 4.364895498681314]"
$ws.Range("R14").Value = 0.428862533655834
$ws.Range("S14").Value = 0.428862533655834
$ws.Range("T14").Value = 11.81662773770112
$ws.Range("U14").Value = "[7.257192068611488, 16.37606340679075]"
$ws.Range("V14").Value = 0.000004408283154511139
$ws.Range("W14").Value = 0.000004408283154511139
$ws.Range("X14").Value = 19.32346346346372
$ws.Range("Y14").Value = 7.348698698698799
$ws.Range("Z14").Value = 31.29822822822865
$ws.Range("F15").Value = 24.07000000000032
$ws.Range("H15").Value = 0.3490703915710098
$ws.Range("I15").Value = 0.3490703915710098
$ws.Range("L15").Value = 4.733846473306594
$ws.Range("M15").Value = "[-4.086434325020129, 13.554127271633316]"
$ws.Range("N15").Value = 0.2854698421333193
$ws.Range("O15").Value = 0.2854698421333193
$ws.Range("P15").Value = 2.182447749340657
$ws.Range("Q15").Value = "[-0.9245527929483481, 5.289448291629662]"
$ws.Range("R15").Value = 0.1640201743115917
$ws.Range("S15").Value = 0.1640201743115917
$ws.Range("T15").Value = 12.05537946472332
$ws.Range("U15").Value = "[7.46669855000529, 16.644060379441342]"
$ws.Range("V15").Value = 0.000003465763093624119
$ws.Range("W15").Value = 0.000003465763093624119
$ws.Range("X15").Value = 15.70934934934956
$ws.Range("Y15").Value = 3.806866866866919
$ws.Range("Z15").Value = 27.61183183183221
